$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Step 2 price + Step1 date edits ---
$ws.Range("A1").Value = 45309
$ws.Range("D20").Value = 565.303
$ws.Range("D21").Value = 626.903
$ws.Range("D22").Value = 808.479
$ws.Range("D23").Value = 1212.716
$ws.Range("D24").Value = 2425.391
$ws.Range("D25").Value = 2911.734
$ws.Range("D26").Value = 486.386
$ws.Range("D27").Value = 587.396
$ws.Range("D28").Value = 546.362
$ws.Range("D29").Value = 647.398
$ws.Range("D30").Value = 546.362
$ws.Range("D31").Value = 647.398

# --- Recreate merged ranges in the order the authoring tool produced, ---
# --- preserving each underlying cell's original style.       ---
$ws.Range("B27:C27").Copy()
$ws.Range("Z200:AA200").PasteSpecial(-4122)
$ws.Range("B27:C27").UnMerge()
$ws.Range("B27:C27").Merge()
$ws.Range("Z200:AA200").Copy()
$ws.Range("B27:C27").PasteSpecial(-4122)
$ws.Range("Z200:AA200").Clear()

$ws.Range("B24:C24").Copy()
$ws.Range("Z201:AA201").PasteSpecial(-4122)
$ws.Range("B24:C24").UnMerge()
$ws.Range("B24:C24").Merge()
$ws.Range("Z201:AA201").Copy()
$ws.Range("B24:C24").PasteSpecial(-4122)
$ws.Range("Z201:AA201").Clear()

$ws.Range("A10:D10").Copy()
$ws.Range("Z202:AC202").PasteSpecial(-4122)
$ws.Range("A10:D10").UnMerge()
$ws.Range("A10:D10").Merge()
$ws.Range("Z202:AC202").Copy()
$ws.Range("A10:D10").PasteSpecial(-4122)
$ws.Range("Z202:AC202").Clear()

$ws.Range("B25:C25").Copy()
$ws.Range("Z203:AA203").PasteSpecial(-4122)
$ws.Range("B25:C25").UnMerge()
$ws.Range("B25:C25").Merge()
$ws.Range("Z203:AA203").Copy()
$ws.Range("B25:C25").PasteSpecial(-4122)
$ws.Range("Z203:AA203").Clear()

$ws.Range("A11:D11").Copy()
$ws.Range("Z204:AC204").PasteSpecial(-4122)
$ws.Range("A11:D11").UnMerge()
$ws.Range("A11:D11").Merge()
$ws.Range("Z204:AC204").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)
$ws.Range("Z204:AC204").Clear()

$ws.Range("B31:C31").Copy()
$ws.Range("Z205:AA205").PasteSpecial(-4122)
$ws.Range("B31:C31").UnMerge()
$ws.Range("B31:C31").Merge()
$ws.Range("Z205:AA205").Copy()
$ws.Range("B31:C31").PasteSpecial(-4122)
$ws.Range("Z205:AA205").Clear()

$ws.Range("B19:C19").Copy()
$ws.Range("Z206:AA206").PasteSpecial(-4122)
$ws.Range("B19:C19").UnMerge()
$ws.Range("B19:C19").Merge()
$ws.Range("Z206:AA206").Copy()
$ws.Range("B19:C19").PasteSpecial(-4122)
$ws.Range("Z206:AA206").Clear()

$ws.Range("A1:D1").Copy()
$ws.Range("Z207:AC207").PasteSpecial(-4122)
$ws.Range("A1:D1").UnMerge()
$ws.Range("A1:D1").Merge()
$ws.Range("Z207:AC207").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)
$ws.Range("Z207:AC207").Clear()

$ws.Range("B21:C21").Copy()
$ws.Range("Z208:AA208").PasteSpecial(-4122)
$ws.Range("B21:C21").UnMerge()
$ws.Range("B21:C21").Merge()
$ws.Range("Z208:AA208").Copy()
$ws.Range("B21:C21").PasteSpecial(-4122)
$ws.Range("Z208:AA208").Clear()

$ws.Range("B20:C20").Copy()
$ws.Range("Z209:AA209").PasteSpecial(-4122)
$ws.Range("B20:C20").UnMerge()
$ws.Range("B20:C20").Merge()
$ws.Range("Z209:AA209").Copy()
$ws.Range("B20:C20").PasteSpecial(-4122)
$ws.Range("Z209:AA209").Clear()

$ws.Range("B30:C30").Copy()
$ws.Range("Z210:AA210").PasteSpecial(-4122)
$ws.Range("B30:C30").UnMerge()
$ws.Range("B30:C30").Merge()
$ws.Range("Z210:AA210").Copy()
$ws.Range("B30:C30").PasteSpecial(-4122)
$ws.Range("Z210:AA210").Clear()

$ws.Range("B26:C26").Copy()
$ws.Range("Z211:AA211").PasteSpecial(-4122)
$ws.Range("B26:C26").UnMerge()
$ws.Range("B26:C26").Merge()
$ws.Range("Z211:AA211").Copy()
$ws.Range("B26:C26").PasteSpecial(-4122)
$ws.Range("Z211:AA211").Clear()

$ws.Range("A9:D9").Copy()
$ws.Range("Z212:AC212").PasteSpecial(-4122)
$ws.Range("A9:D9").UnMerge()
$ws.Range("A9:D9").Merge()
$ws.Range("Z212:AC212").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)
$ws.Range("Z212:AC212").Clear()

$ws.Range("B22:C22").Copy()
$ws.Range("Z213:AA213").PasteSpecial(-4122)
$ws.Range("B22:C22").UnMerge()
$ws.Range("B22:C22").Merge()
$ws.Range("Z213:AA213").Copy()
$ws.Range("B22:C22").PasteSpecial(-4122)
$ws.Range("Z213:AA213").Clear()

$ws.Range("B28:C28").Copy()
$ws.Range("Z214:AA214").PasteSpecial(-4122)
$ws.Range("B28:C28").UnMerge()
$ws.Range("B28:C28").Merge()
$ws.Range("Z214:AA214").Copy()
$ws.Range("B28:C28").PasteSpecial(-4122)
$ws.Range("Z214:AA214").Clear()

$ws.Range("B23:C23").Copy()
$ws.Range("Z215:AA215").PasteSpecial(-4122)
$ws.Range("B23:C23").UnMerge()
$ws.Range("B23:C23").Merge()
$ws.Range("Z215:AA215").Copy()
$ws.Range("B23:C23").PasteSpecial(-4122)
$ws.Range("Z215:AA215").Clear()

$ws.Range("B29:C29").Copy()
$ws.Range("Z216:AA216").PasteSpecial(-4122)
$ws.Range("B29:C29").UnMerge()
$ws.Range("B29:C29").Merge()
$ws.Range("Z216:AA216").Copy()
$ws.Range("B29:C29").PasteSpecial(-4122)
$ws.Range("Z216:AA216").Clear()
